$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D (Fecha), K (Precio mínimo), L (Precio máximo),
# M (Precio promedio ponderado), O (Origen), P (Precio $/Kg) for rows 2-11.

$rows = @(
    @{ Row = 2;  D = 44162; K = 7500;  L = 8000;  M = 7750;  O = "Región Metropolitana"; P = 310 },
    @{ Row = 3;  D = 44160; K = 9000;  L = 10000; M = 9500;  O = "Región Metropolitana"; P = 380 },
    @{ Row = 4;  D = 44188; K = 18000; L = 20000; M = 19000; O = "Región Metropolitana"; P = 760 },
    @{ Row = 5;  D = 44316; K = 16000; L = 18000; M = 17000; O = "Región Metropolitana"; P = 680 },
    @{ Row = 6;  D = 44454; K = 13000; L = 14000; M = 13500; O = "Provincia del Elquí";  P = 540 },
    @{ Row = 7;  D = 44467; K = 8000;  L = 9000;  M = 8500;  O = "Región Metropolitana"; P = 340 },
    @{ Row = 8;  D = 44351; K = 15000; L = 16000; M = 15500; O = "Región Metropolitana"; P = 620 },
    @{ Row = 9;  D = 44335; K = 18000; L = 20000; M = 19000; O = "Provincia de Limarí";  P = 760 },
    @{ Row = 10; D = 44384; K = 12000; L = 13000; M = 12500; O = "Región de Coquimbo";   P = 500 },
    @{ Row = 11; D = 44461; K = 13000; L = 14000; M = 13500; O = "Provincia del Elquí";  P = 540 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio mínimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio máximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
}
